$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Junio de 2020 a las 21:30"

# Country-name swaps caused by shared-string reordering (same rows, swapped display text)
$ws.Range("A94").Value = "Mauritania"
$ws.Range("A95").Value = "Grecia"
$ws.Range("A100").Value = "Costa Rica"
$ws.Range("A101").Value = "Mayotte"
$ws.Range("A102").Value = "Croacia"
$ws.Range("A119").Value = "Paraguay"
$ws.Range("A120").Value = "Nueva Zelanda"
$ws.Range("A121").Value = "Zambia"
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"

# Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 2447449
$ws.Range("C4").Value = 23281
$ws.Range("D4").Value = 1025350
$ws.Range("E4").Value = 1298165
$ws.Range("G4").Value = 461
$ws.Range("H4").Value = 123934
$ws.Range("B7").Value = 472972
$ws.Range("C7").Value = 16857
$ws.Range("D7").Value = 271688
$ws.Range("E7").Value = 186377
$ws.Range("G7").Value = 424
$ws.Range("H7").Value = 14907
$ws.Range("D11").Value = 215093
$ws.Range("E11").Value = 34592
$ws.Range("B14").Value = 192977
$ws.Range("C14").Value = 199
$ws.Range("E14").Value = 7687
$ws.Range("G14").Value = 4
$ws.Range("H14").Value = 8990
$ws.Range("B19").Value = 161348
$ws.Range("C19").Value = 81
$ws.Range("D19").Value = 75127
$ws.Range("E19").Value = 56490
$ws.Range("G19").Value = 11
$ws.Range("H19").Value = 29731
$ws.Range("B58").Value = 15013
$ws.Range("C58").Value = 445
$ws.Range("D58").Value = 11078
$ws.Range("E58").Value = 3840
$ws.Range("B76").Value = 6901
$ws.Range("C76").Value = 239
$ws.Range("E76").Value = 2197
$ws.Range("B94").Value = 3519
$ws.Range("C94").Value = 227
$ws.Range("D94").Value = 1074
$ws.Range("E94").Value = 2329
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 116
$ws.Range("B95").Value = 3310
$ws.Range("C95").Value = 8
$ws.Range("D95").Value = 1374
$ws.Range("E95").Value = 1746
$ws.Range("H95").Value = 190
$ws.Range("B100").Value = 2515
$ws.Range("C100").Value = 147
$ws.Range("D100").Value = 1210
$ws.Range("E100").Value = 1293
$ws.Range("H100").Value = 12
$ws.Range("B101").Value = 2467
$ws.Range("C101").Value = 33
$ws.Range("D101").Value = 2218
$ws.Range("E101").Value = 217
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 32
$ws.Range("B102").Value = 2388
$ws.Range("C102").Value = 22
$ws.Range("D102").Value = 2145
$ws.Range("E102").Value = 136
$ws.Range("H102").Value = 107
$ws.Range("B119").Value = 1528
$ws.Range("C119").Value = 106
$ws.Range("D119").Value = 944
$ws.Range("E119").Value = 571
$ws.Range("H119").Value = 13
$ws.Range("B120").Value = 1516
$ws.Range("C120").Value = 1
$ws.Range("D120").Value = 1483
$ws.Range("E120").Value = 11
$ws.Range("H120").Value = 22
$ws.Range("B121").Value = 1489
$ws.Range("C121").Value = 12
$ws.Range("D121").Value = 1223
$ws.Range("E121").Value = 248
$ws.Range("H121").Value = 18
$ws.Range("D137").Value = 770
$ws.Range("E137").Value = 16
$ws.Range("B146").Value = 690
$ws.Range("C146").Value = 16
$ws.Range("D146").Value = 331
$ws.Range("E146").Value = 352
$ws.Range("B149").Value = 662
$ws.Range("C149").Value = 10
$ws.Range("E149").Value = 358
$ws.Range("B154").Value = 508
$ws.Range("C154").Value = 1
$ws.Range("E154").Value = 47
$ws.Range("D160").Value = 142
$ws.Range("E160").Value = 168
